$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 2 new rows before the current row 235 (Lutjanidae 12S rows) ---
$ws.Rows("235:236").Insert()

# --- Insert 3 new rows before the current row 242 (now between Lutjanidae ---
# --- 16S "Genus and species wrong" row and the Rottnest 12S block) ---
$ws.Rows("242:244").Insert()

# --- Fill the 2 newly inserted rows at 235:236 ("100 Australian species" / CO1) ---
$ws.Range("A235").Value = "VSEARCH"
$ws.Range("B235").Value = "100 Australian species"
$ws.Range("C235").Value = "CO1"
$ws.Range("D235").Value = "Genus and species correct"
$ws.Range("E235").Value = 1
$ws.Range("F235").Value = 5
$ws.Range("G235").Value = 20

$ws.Range("A236").Value = "VSEARCH"
$ws.Range("B236").Value = "100 Australian species"
$ws.Range("C236").Value = "CO1"
$ws.Range("D236").Value = "Genus correct, species wrong"
$ws.Range("E236").Value = 4
$ws.Range("F236").Value = 5
$ws.Range("G236").Value = 80

# --- Fill the 3 newly inserted rows at 242:244 (Lutjanidae/CO1 + Rottnest/CO1) ---
$ws.Range("A242").Value = "VSEARCH"
$ws.Range("B242").Value = "Lutjanidae"
$ws.Range("C242").Value = "CO1"
$ws.Range("D242").Value = "Genus and species correct"
$ws.Range("E242").Value = 5
$ws.Range("F242").Value = 5
$ws.Range("G242").Value = 100

$ws.Range("A243").Value = "VSEARCH"
$ws.Range("B243").Value = "Rottnest"
$ws.Range("C243").Value = "CO1"
$ws.Range("D243").Value = "Genus and species correct"
$ws.Range("E243").Value = 25
$ws.Range("F243").Value = 34
$ws.Range("G243").Value = 73.53

$ws.Range("A244").Value = "VSEARCH"
$ws.Range("B244").Value = "Rottnest"
$ws.Range("C244").Value = "CO1"
$ws.Range("D244").Value = "Genus correct, species wrong"
$ws.Range("E244").Value = 9
$ws.Range("F244").Value = 34
$ws.Range("G244").Value = 26.47
